$wb = $excel.ActiveWorkbook

# This is a "Generate Report for Handoff" refresh: the 84377e8a, 9c2f69b4,
# ef67b2d3 and f4c95d07 files (rows 4-7 on each sheet) were just handed off
# again, so their priority flips from "low" to "ht" and the handoff
# timestamps for that batch move forward.

# "Overview" sheet: Latest HO Xliff Generate Date for rows 4-7 moves from
# 2016-09-01 12:34:15 to 2016-09-01 12:34:38.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-09-01 12:34:38"

# "zh-cn" sheet: Priority for rows 4-7 goes from "low" to "ht", and the
# Latest Handoff Datetime advances from 2016-09-01 12:34:01 to
# 2016-09-01 12:34:33.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4:E7").Value = "ht"
$wsZh.Range("H4:H7").Value = "2016-09-01 12:34:33"

# "de-de" sheet: same Priority change for rows 4-7, and the Latest Handoff
# Datetime (shared string with the Overview generate date) advances from
# 2016-09-01 12:34:15 to 2016-09-01 12:34:38.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4:E7").Value = "ht"
$wsDe.Range("H4:H7").Value = "2016-09-01 12:34:38"
